$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 245.5
$ws.Range("I5").Value = 245.5
$ws.Range("K5").Value = 245.5
$ws.Range("M5").Value = -130.5
# Row 6
$ws.Range("H6").Value = 94.454544
$ws.Range("I6").Value = 116.625
$ws.Range("J6").Value = 35.333332
$ws.Range("K6").Value = 349.875
$ws.Range("L6").Value = 105.999996
$ws.Range("M6").Value = -237.875
$ws.Range("N6").Value = -329.999996
# Row 8
$ws.Range("H8").Value = 10
$ws.Range("I8").Value = 10
$ws.Range("K8").Value = 30
$ws.Range("M8").Value = 109
# Row 12
$ws.Range("H12").Value = 498.30768
$ws.Range("I12").Value = 481.5
$ws.Range("J12").Value = 700
$ws.Range("K12").Value = 481.5
$ws.Range("L12").Value = 700
$ws.Range("M12").Value = -311.5
$ws.Range("N12").Value = -1040
# Row 13
$ws.Range("H13").Value = 1000
$ws.Range("J13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("N13").Value = -1338
# Row 21
$ws.Range("H21").Value = 22500
$ws.Range("I21").Value = 15000
$ws.Range("K21").Value = 15000
$ws.Range("M21").Value = -14532
# Row 23
$ws.Range("H23").Value = 22500
$ws.Range("I23").Value = 15000
$ws.Range("K23").Value = 15000
$ws.Range("M23").Value = -14766
# Row 92
$ws.Range("H92").Value = 125000870
$ws.Range("I92").Value = 200000750
$ws.Range("J92").Value = 1066.6666
$ws.Range("K92").Value = 200000750
$ws.Range("L92").Value = 1066.6666
$ws.Range("M92").Value = -199999502
$ws.Range("N92").Value = -3562.6666
# Row 137
$ws.Range("H137").Value = 1826.3334
$ws.Range("I137").Value = 1250
$ws.Range("K137").Value = 3750
$ws.Range("M137").Value = -1200
# Row 138
$ws.Range("H138").Value = 7483.6924
$ws.Range("J138").Value = 7208.8
$ws.Range("L138").Value = 21626.4
$ws.Range("N138").Value = -31906.4

$ws = $wb.Worksheets.Item("ARM")
# Row 46
$ws.Range("H46").Value = 19033.75
$ws.Range("I46").Value = 19033.75
$ws.Range("K46").Value = 19033.75
$ws.Range("M46").Value = -18714.75
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
# Row 97
$ws.Range("H97").Value = 37038840
$ws.Range("I97").Value = 37038840
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 37038840
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -37038344
$ws.Range("N97").ClearContents()
# Row 102
$ws.Range("H102").Value = 30002656
$ws.Range("I102").Value = 1669764.1
$ws.Range("K102").Value = 1669764.1
$ws.Range("M102").Value = -1668142.1
# Row 115
$ws.Range("H115").Value = 41500
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 41500
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 41500
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -44634

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 140119.12
$ws.Range("I94").Value = 185701.33
$ws.Range("J94").Value = 3372.5
$ws.Range("K94").Value = 185701.33
$ws.Range("L94").Value = 3372.5
$ws.Range("M94").Value = -185250.33
$ws.Range("N94").Value = -4274.5
# Row 132
$ws.Range("H132").Value = 60780
$ws.Range("J132").Value = 60780
$ws.Range("L132").Value = 60780
$ws.Range("N132").Value = -70900

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 197.90909
$ws.Range("I7").Value = 37.4
$ws.Range("J7").Value = 331.66666
$ws.Range("K7").Value = 37.4
$ws.Range("L7").Value = 331.66666
$ws.Range("M7").Value = 75.59999999999999
$ws.Range("N7").Value = -557.66666
# Row 16
$ws.Range("H16").Value = 1919.2858
$ws.Range("J16").Value = 2219.4
$ws.Range("L16").Value = 2219.4
$ws.Range("N16").Value = -2793.4
# Row 18
$ws.Range("H18").Value = 77000
$ws.Range("J18").Value = 77000
$ws.Range("L18").Value = 77000
$ws.Range("N18").Value = -77460
# Row 25
$ws.Range("H25").Value = 25000
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
# Row 113
$ws.Range("H113").Value = 1919.2858
$ws.Range("J113").Value = 2219.4
$ws.Range("L113").Value = 2219.4
$ws.Range("N113").Value = -6559.4

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 198.2
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 198.2
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 594.5999999999999
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -940.5999999999999
# Row 59
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
# Row 68
$ws.Range("H68").Value = 985.7778
$ws.Range("J68").Value = 660
$ws.Range("L68").Value = 1980
$ws.Range("N68").Value = -3602
# Row 71
$ws.Range("H71").Value = 985.7778
$ws.Range("J71").Value = 660
$ws.Range("L71").Value = 5940
$ws.Range("N71").Value = -14052

$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 8007.909
$ws.Range("J24").Value = 8007.909
$ws.Range("L24").Value = 8007.909
$ws.Range("N24").Value = -8353.909
# Row 97
$ws.Range("H97").Value = 4958.3335
$ws.Range("I97").Value = 4950
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 4950
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -4454
$ws.Range("N97").Value = -5992
# Row 122
$ws.Range("H122").Value = 2127.2354
$ws.Range("I122").Value = 2368.7
$ws.Range("J122").Value = 1782.2858
$ws.Range("K122").Value = 7106.099999999999
$ws.Range("L122").Value = 5346.857400000001
$ws.Range("M122").Value = -4656.099999999999
$ws.Range("N122").Value = -10246.8574

$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 3000
$ws.Range("J20").Value = 3000
$ws.Range("L20").Value = 3000
$ws.Range("N20").Value = -3452
# Row 40
$ws.Range("H40").Value = 2066
$ws.Range("I40").Value = 1100
$ws.Range("K40").Value = 1100
$ws.Range("M40").Value = -964

$ws = $wb.Worksheets.Item("WVR")
# Row 22
$ws.Range("H22").Value = 999
$ws.Range("J22").Value = 999
$ws.Range("L22").Value = 999
$ws.Range("N22").Value = -1585
# Row 100
$ws.Range("H100").Value = 4271.8335
$ws.Range("I100").Value = 3729.3333
$ws.Range("J100").Value = 4814.3335
$ws.Range("K100").Value = 7458.6666
$ws.Range("L100").Value = 9628.666999999999
$ws.Range("M100").Value = -6917.6666
$ws.Range("N100").Value = -10710.667
# Row 130
$ws.Range("H130").Value = 40429
$ws.Range("J130").Value = 40429
$ws.Range("L130").Value = 40429
$ws.Range("N130").Value = -50469
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

Write-Output "Done applying Golem Profits updates"
